$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 / 38: Mantle <-> Filecoin swap (name, link, price, volume) ---
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.13'
$ws.Range("E37").Value = '  +3.16%  '

$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.04'
$ws.Range("E38").Value = '  +0.23%  '

# --- Remaining Price / Volume(1h) updates ---
$ws.Range("D2").Value = '69.811.10'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '3.933.29'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.26'
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.33'
$ws.Range("E6").Value = '  +2.60%  '
$ws.Range("D7").Value = '3.935.76'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").Value = '  +4.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.00'
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").Value = '4.592.65'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '3.923.40'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '69.820.34'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.54'
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.06'
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '499.03'
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  +5.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.43'
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '4.085.26'
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.86'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.21'
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("D35").Value = '3.903.96'
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.140'
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.30'
$ws.Range("E40").Value = '  +8.99%  '
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.05'
$ws.Range("E43").Value = '  +2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '435.32'
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.29'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("E48").Value = '  +22.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0367'
$ws.Range("E49").Value = '  +2.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.34'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.96'
$ws.Range("E51").Value = '  +0.90%  '
